$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "venue"
$ws.Cells.Item(1,2).Value = "date"
$ws.Cells.Item(1,3).Value = "result"
$ws.Cells.Item(1,4).Value = "ownTeam"
$ws.Cells.Item(1,5).Value = "oppTeam"
$ws.Cells.Item(1,6).Value = "batsman"
$ws.Cells.Item(1,7).Value = "totalRuns"
$ws.Cells.Item(1,8).Value = "totalBalls"
$ws.Cells.Item(1,9).Value = "total4s"
$ws.Cells.Item(1,10).Value = "total6s"
$ws.Cells.Item(1,11).Value = "sr"

$ws.Cells.Item(2,1).Value = " Dubai (DSC)"
$ws.Cells.Item(2,2).Value = " October 17 2020"
$ws.Cells.Item(2,3).Value = "RCB won by 7 wickets (with 2 balls remaining)"
$ws.Cells.Item(2,4).Value = "Rajasthan Royals"
$ws.Cells.Item(2,5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(2,6).Value = "Jofra Archer "
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = "2"
$ws.Cells.Item(2,8).NumberFormat = "@"
$ws.Cells.Item(2,8).Value = "3"
$ws.Cells.Item(2,9).NumberFormat = "@"
$ws.Cells.Item(2,9).Value = "0"
$ws.Cells.Item(2,10).NumberFormat = "@"
$ws.Cells.Item(2,10).Value = "0"
$ws.Cells.Item(2,11).NumberFormat = "@"
$ws.Cells.Item(2,11).Value = "66.66"

$ws.Cells.Item(3,1).Value = " Sharjah"
$ws.Cells.Item(3,2).Value = " September 27 2020"
$ws.Cells.Item(3,3).Value = "Royals won by 4 wickets (with 3 balls remaining)"
$ws.Cells.Item(3,4).Value = "Rajasthan Royals"
$ws.Cells.Item(3,5).Value = "Kings XI Punjab"
$ws.Cells.Item(3,6).Value = "Jofra Archer "
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = "13"
$ws.Cells.Item(3,8).NumberFormat = "@"
$ws.Cells.Item(3,8).Value = "3"
$ws.Cells.Item(3,9).NumberFormat = "@"
$ws.Cells.Item(3,9).Value = "0"
$ws.Cells.Item(3,10).NumberFormat = "@"
$ws.Cells.Item(3,10).Value = "2"
$ws.Cells.Item(3,11).NumberFormat = "@"
$ws.Cells.Item(3,11).Value = "433.33"

$ws.Cells.Item(4,1).Value = " Dubai (DSC)"
$ws.Cells.Item(4,2).Value = " September 30 2020"
$ws.Cells.Item(4,3).Value = "KKR won by 37 runs"
$ws.Cells.Item(4,4).Value = "Rajasthan Royals"
$ws.Cells.Item(4,5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(4,6).Value = "Jofra Archer "
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = "6"
$ws.Cells.Item(4,8).NumberFormat = "@"
$ws.Cells.Item(4,8).Value = "4"
$ws.Cells.Item(4,9).NumberFormat = "@"
$ws.Cells.Item(4,9).Value = "0"
$ws.Cells.Item(4,10).NumberFormat = "@"
$ws.Cells.Item(4,10).Value = "1"
$ws.Cells.Item(4,11).NumberFormat = "@"
$ws.Cells.Item(4,11).Value = "150.00"

$ws.Cells.Item(5,1).Value = " Abu Dhabi"
$ws.Cells.Item(5,2).Value = " October 03 2020"
$ws.Cells.Item(5,3).Value = "RCB won by 8 wickets (with 5 balls remaining)"
$ws.Cells.Item(5,4).Value = "Rajasthan Royals"
$ws.Cells.Item(5,5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(5,6).Value = "Jofra Archer "
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = "16"
$ws.Cells.Item(5,8).NumberFormat = "@"
$ws.Cells.Item(5,8).Value = "10"
$ws.Cells.Item(5,9).NumberFormat = "@"
$ws.Cells.Item(5,9).Value = "1"
$ws.Cells.Item(5,10).NumberFormat = "@"
$ws.Cells.Item(5,10).Value = "1"
$ws.Cells.Item(5,11).NumberFormat = "@"
$ws.Cells.Item(5,11).Value = "160.00"

$ws.Cells.Item(6,1).Value = " Dubai (DSC)"
$ws.Cells.Item(6,2).Value = " October 14 2020"
$ws.Cells.Item(6,3).Value = "Capitals won by 13 runs"
$ws.Cells.Item(6,4).Value = "Rajasthan Royals"
$ws.Cells.Item(6,5).Value = "Delhi Capitals"
$ws.Cells.Item(6,6).Value = "Jofra Archer "
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = "1"
$ws.Cells.Item(6,8).NumberFormat = "@"
$ws.Cells.Item(6,8).Value = "4"
$ws.Cells.Item(6,9).NumberFormat = "@"
$ws.Cells.Item(6,9).Value = "0"
$ws.Cells.Item(6,10).NumberFormat = "@"
$ws.Cells.Item(6,10).Value = "0"
$ws.Cells.Item(6,11).NumberFormat = "@"
$ws.Cells.Item(6,11).Value = "25.00"

$ws.Cells.Item(7,1).Value = " Sharjah"
$ws.Cells.Item(7,2).Value = " October 09 2020"
$ws.Cells.Item(7,3).Value = "Capitals won by 46 runs"
$ws.Cells.Item(7,4).Value = "Rajasthan Royals"
$ws.Cells.Item(7,5).Value = "Delhi Capitals"
$ws.Cells.Item(7,6).Value = "Jofra Archer "
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = "2"
$ws.Cells.Item(7,8).NumberFormat = "@"
$ws.Cells.Item(7,8).Value = "4"
$ws.Cells.Item(7,9).NumberFormat = "@"
$ws.Cells.Item(7,9).Value = "0"
$ws.Cells.Item(7,10).NumberFormat = "@"
$ws.Cells.Item(7,10).Value = "0"
$ws.Cells.Item(7,11).NumberFormat = "@"
$ws.Cells.Item(7,11).Value = "50.00"

$ws.Cells.Item(8,1).Value = " Dubai (DSC)"
$ws.Cells.Item(8,2).Value = " October 22 2020"
$ws.Cells.Item(8,3).Value = "Sunrisers won by 8 wickets (with 11 balls remaining)"
$ws.Cells.Item(8,4).Value = "Rajasthan Royals"
$ws.Cells.Item(8,5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(8,6).Value = "Jofra Archer "
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = "16"
$ws.Cells.Item(8,8).NumberFormat = "@"
$ws.Cells.Item(8,8).Value = "7"
$ws.Cells.Item(8,9).NumberFormat = "@"
$ws.Cells.Item(8,9).Value = "1"
$ws.Cells.Item(8,10).NumberFormat = "@"
$ws.Cells.Item(8,10).Value = "1"
$ws.Cells.Item(8,11).NumberFormat = "@"
$ws.Cells.Item(8,11).Value = "228.57"

$ws.Cells.Item(9,1).Value = " Dubai (DSC)"
$ws.Cells.Item(9,2).Value = " November 01 2020"
$ws.Cells.Item(9,3).Value = "KKR won by 60 runs"
$ws.Cells.Item(9,4).Value = "Rajasthan Royals"
$ws.Cells.Item(9,5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(9,6).Value = "Jofra Archer "
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = "6"
$ws.Cells.Item(9,8).NumberFormat = "@"
$ws.Cells.Item(9,8).Value = "9"
$ws.Cells.Item(9,9).NumberFormat = "@"
$ws.Cells.Item(9,9).Value = "0"
$ws.Cells.Item(9,10).NumberFormat = "@"
$ws.Cells.Item(9,10).Value = "0"
$ws.Cells.Item(9,11).NumberFormat = "@"
$ws.Cells.Item(9,11).Value = "66.66"

$ws.Cells.Item(10,1).Value = " Sharjah"
$ws.Cells.Item(10,2).Value = " September 22 2020"
$ws.Cells.Item(10,3).Value = "Royals won by 16 runs"
$ws.Cells.Item(10,4).Value = "Rajasthan Royals"
$ws.Cells.Item(10,5).Value = "Chennai Super Kings"
$ws.Cells.Item(10,6).Value = "Jofra Archer "
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = "27"
$ws.Cells.Item(10,8).NumberFormat = "@"
$ws.Cells.Item(10,8).Value = "8"
$ws.Cells.Item(10,9).NumberFormat = "@"
$ws.Cells.Item(10,9).Value = "0"
$ws.Cells.Item(10,10).NumberFormat = "@"
$ws.Cells.Item(10,10).Value = "4"
$ws.Cells.Item(10,11).NumberFormat = "@"
$ws.Cells.Item(10,11).Value = "337.50"

$ws.Cells.Item(11,1).Value = " Abu Dhabi"
$ws.Cells.Item(11,2).Value = " October 06 2020"
$ws.Cells.Item(11,3).Value = "Mumbai won by 57 runs"
$ws.Cells.Item(11,4).Value = "Rajasthan Royals"
$ws.Cells.Item(11,5).Value = "Mumbai Indians"
$ws.Cells.Item(11,6).Value = "Jofra Archer "
$ws.Cells.Item(11,7).NumberFormat = "@"
$ws.Cells.Item(11,7).Value = "24"
$ws.Cells.Item(11,8).NumberFormat = "@"
$ws.Cells.Item(11,8).Value = "11"
$ws.Cells.Item(11,9).NumberFormat = "@"
$ws.Cells.Item(11,9).Value = "3"
$ws.Cells.Item(11,10).NumberFormat = "@"
$ws.Cells.Item(11,10).Value = "1"
$ws.Cells.Item(11,11).NumberFormat = "@"
$ws.Cells.Item(11,11).Value = "218.18"
